# Insert a new weekly price record as row 88 (pushing the existing rows
# 88-172 down to 89-173) in the daily-logic subset sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 88..172 down to 89..173, leaving a blank row 88 to fill in.
$ws.Rows.Item(88).Insert()

# Populate the new row 88 with the new record's data.
$ws.Range("A88").Value = 4
$ws.Range("B88").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C88").Value = "Los Lagos"
$ws.Range("D88").Value = 44566
$ws.Range("E88").Value = 10
$ws.Range("F88").Value = 100112032
$ws.Range("G88").Value = "Zapallo italiano"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 35
$ws.Range("K88").Value = 11000
$ws.Range("L88").Value = 11000
$ws.Range("M88").Value = 11000
$ws.Range("N88").Value = '$/caja 50 unidades'
$ws.Range("O88").Value = "Región de O'Higgins"
$ws.Range("P88").Value = 220
$ws.Range("Q88").Value = 50
$ws.Range("R88").Value = "Hortaliza"
